$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames (row 1) ---
# K1: seas_id_x -> seas_id
$ws.Range("K1").Value() = "seas_id"
# L1: player_id -> player_id_x
$ws.Range("L1").Value() = "player_id_x"
# N1: seas_id_y -> season_ending_year_y
$ws.Range("N1").Value() = "season_ending_year_y"
# O1: season_ending_year_y -> player_id_y
$ws.Range("O1").Value() = "player_id_y"

# --- Row 2 data: N2/O2 swap meaning ---
# N2 was the numeric seas_id_y (31786); it now holds the season_ending_year_y
# text value "2024" (same text already used elsewhere in the sheet).
$ws.Range("N2").Value() = "'2024"
$ws.Range("N2").Style = "Normal"
# O2 was the text season_ending_year_y ("2024"); it now holds the numeric
# player_id_y value.
$ws.Range("O2").Value() = 4666

# --- Row 3 data: N3/O3 swap meaning ---
$ws.Range("N3").Value() = "'2023"
$ws.Range("N3").Style = "Normal"
$ws.Range("O3").Value() = 1255
